{"js": "// Highlight (yellow) the paragraph mark and run text of ten specific\n// checklist items in the \"Build An Image Gallery\" assignment. These are\n// matched by their exact paragraph text so the edit is robust to any\n// reordering of the document.\nconst targetTexts = [\n  \"The \\u201cgallery\\u201d id.\",\n  \"This will go inside the main tag.\",\n  \"Set this tag to be a flex layout.\",\n  \"Make sure this flex element will wrap the images.\",\n  \"Try to get the flex element to display three images per row of images.\",\n  \"Use a \\u201cjustify-content\\u201d line to make the content look nice.\",\n  \"New selector! Look up \\u201c:hover\\u201d on w3schools. 10 points.\",\n  \"When the user hovers over an image with their mouse, rotate the image 5 degrees and give it a pleasing border color.\",\n  \"The footer. Five points.\",\n  \"Do I need to say it? Make it look nice.\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  paragraph.load(\"text\");\n}\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (targetTexts.indexOf(paragraph.text) !== -1) {\n    // Setting highlightColor on the paragraph's font applies the\n    // <w:highlight> both to the paragraph mark run properties (pPr/rPr)\n    // and to the run(s) that make up the paragraph's text.\n    paragraph.font.highlightColor = \"Yellow\";\n  }\n}\nawait context.sync();\n", "ps1": "# Highlight (yellow) ten specific checklist items in the \"Build An Image\n# Gallery\" assignment. Each paragraph is matched by its exact text so the\n# edit is robust to any reordering of the document. Setting\n# HighlightColorIndex on the paragraph range's Font (rather than on the\n# Range itself) applies <w:highlight> to both the paragraph mark run\n# properties (pPr/rPr) and the run(s) that hold the paragraph's text,\n# matching how Word highlights an entire paragraph (including its pilcrow).\n$d = $word.ActiveDocument\n\n$targetTexts = @(\n    \"The \u201cgallery\u201d id.\",\n    \"This will go inside the main tag.\",\n    \"Set this tag to be a flex layout.\",\n    \"Make sure this flex element will wrap the images.\",\n    \"Try to get the flex element to display three images per row of images.\",\n    \"Use a \u201cjustify-content\u201d line to make the content look nice.\",\n    \"New selector! Look up \u201c:hover\u201d on w3schools. 10 points.\",\n    \"When the user hovers over an image with their mouse, rotate the image 5 degrees and give it a pleasing border color.\",\n    \"The footer. Five points.\",\n    \"Do I need to say it? Make it look nice.\"\n)\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($targetTexts -contains $text) {\n        $p.Range.Font.HighlightColorIndex = \"Yellow\"\n    }\n}\n"}
